{"js": "// Insert a comma after \"En esta secci\u00f3n\" so the sentence reads\n// \"En esta secci\u00f3n, se explorar\u00e1 la estructura inicial del proyecto...\".\nconst body = context.document.body;\n\n// The phrase is unique in the document, so a single search hit is expected.\nconst results = body.search(\"En esta secci\u00f3n\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items/text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the target sentence start 'En esta secci\u00f3n'.\");\n}\n\n// Insert the comma immediately after the matched phrase (before the space\n// that precedes \"se explorar\u00e1...\").\nconst hit = results.items[0];\nhit.insertText(\",\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Insert a comma after \"En esta secci\u00f3n\" so the sentence reads\n# \"En esta secci\u00f3n, se explorar\u00e1 la estructura inicial del proyecto...\".\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$find = $range.Find\n$find.ClearFormatting()\n$find.Text = \"secci\u00f3n\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n\n$found = $find.Execute()\nif (-not $found) {\n    throw \"Could not find the target word 'secci\u00f3n'.\"\n}\n\n# $range now spans the matched word (\"secci\u00f3n\"); collapse it to its end\n# point (0 = wdCollapseEnd) and insert the comma right there.\n$range.Collapse(0)\n$range.InsertAfter(\",\")\n"}
